# edit.ps1 - Apply "Update countries & provincias Spain" changes
# Reorders several country rows (data refreshed / re-sorted) and updates
# the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados ..." timestamp banner (A1): 13:22 -> 13:52
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 13:52"

# Table updates: Row, new country name (or $null to keep as-is),
# then Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes (columns B..H)
$updates = @(
    @{ Row = 15; Name = $null; Vals = @(18803, 952, 250, 16686, 1385, 101, 1867) },
    @{ Row = 70; Name = "Barein"; Vals = @(723, 23, 451, 268, 4, 0, 4) },
    @{ Row = 71; Name = "Crucero"; Vals = @(712, 0, 619, 82, 10, 0, 11) },
    @{ Row = 72; Name = "Bielorrusia"; Vals = @(700, 138, 53, 634, 11, 5, 13) },
    @{ Row = 73; Name = $null; Vals = @(667, 13, 47, 593, 4, 4, 27) },
    @{ Row = 76; Name = "Azerbaiyan"; Vals = @(641, 57, 44, 590, 17, 0, 7) },
    @{ Row = 77; Name = "Kazajistan"; Vals = @(629, 45, 45, 578, 16, 0, 6) },
    @{ Row = 101; Name = "Vietnam"; Vals = @(245, 4, 95, 150, 8, 0, 0) },
    @{ Row = 102; Name = "Malta"; Vals = @(241, 14, 5, 236, 3, 0, 0) },
    @{ Row = 158; Name = $null; Vals = @(22, 0, 5, 16, 0, 1, 1) },
    @{ Row = 207; Name = "Burundi"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 208; Name = "Islas Virgenes Britanicas"; Vals = @(3, 0, 0, 3, 0, 0, 0) },
    @{ Row = 209; Name = "Islas Malvinas"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 210; Name = "Bonaire, San Eustaquio y Saba"; Vals = @(2, 0, 0, 2, 0, 0, 0) },
    @{ Row = 211; Name = "Papua Nueva Guinea"; Vals = @(2, 1, 0, 2, 0, 0, 0) },
    @{ Row = 212; Name = "Sudan del Sur"; Vals = @(1, 0, 0, 1, 0, 0, 0) },
    @{ Row = 213; Name = "Timor Oriental"; Vals = @(1, 0, 0, 1, 0, 0, 0) }
)

foreach ($u in $updates) {
    if ($u.Name -ne $null) {
        $ws.Cells.Item($u.Row, 1).Value = $u.Name
    }
    $col = 2
    foreach ($v in $u.Vals) {
        $ws.Cells.Item($u.Row, $col).Value = $v
        $col = $col + 1
    }
}
